# Update "Horarios Línea 141" workbook with the latest scraped schedule data.
# New scrape timestamp: 01:41:13 (previously 00:59:16)

$wb = $excel.ActiveWorkbook

$newTime = "01:41:13"

# --- Sheet "LP1912" (sheet1) ---
# Rows shrink from 3 data rows (6-8) to 2 data rows (6-7); values refreshed.
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 2"

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "01:58"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 17
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "02:58"
$ws1.Range("C7").Value = "215_ALUAR"
$ws1.Range("D7").Value = 77
$ws1.Range("E7").Value = "LP1912"

# Remove the now-unused former row 8
$ws1.Rows.Item(8).Delete()

# --- Sheet "LP1912-215" (sheet2) ---
# Rows shrink from 2 data rows (6-7) to 1 data row (6); values refreshed.
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "02:58"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 77
$ws2.Range("E6").Value = "LP1912"

# Remove the now-unused former row 7
$ws2.Rows.Item(7).Delete()

# --- Sheet "6203-6173" (sheet3) ---
# Only the "Última actualización" timestamp changes; still zero data rows.
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
